$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-9 (columns B:F) with new values
$ws.Range("B2").Value = 0.4769321400286752
$ws.Range("C2").Value = 0.9320513987497278
$ws.Range("D2").Value = 1.650466364378836
$ws.Range("E2").Value = 1.284704777129297
$ws.Range("F2").Value = 1.237927218396544

$ws.Range("B3").Value = 0.6221863334236316
$ws.Range("C3").Value = 0.9661009109131644
$ws.Range("D3").Value = 1.561161447007451
$ws.Range("E3").Value = 1.249464464083493
$ws.Range("F3").Value = 1.127777791928594

$ws.Range("B4").Value = 0.7378115908420845
$ws.Range("C4").Value = 0.9692599323336366
$ws.Range("D4").Value = 1.482360141223014
$ws.Range("E4").Value = 1.217522131717947
$ws.Range("F4").Value = 1.011566308962359

$ws.Range("B5").Value = 0.7806678668064955
$ws.Range("C5").Value = 0.9619392551853039
$ws.Range("D5").Value = 2.031973731581155
$ws.Range("E5").Value = 1.425473160596563
$ws.Range("F5").Value = 1.250913488075273

$ws.Range("B6").Value = 0.7978750242939276
$ws.Range("C6").Value = 0.8777464132826209
$ws.Range("D6").Value = 1.683130036913681
$ws.Range("E6").Value = 1.297355015758478
$ws.Range("F6").Value = 1.078334869922473

$ws.Range("B7").Value = 0.8607107567659571
$ws.Range("C7").Value = 0.8607107567659571
$ws.Range("D7").Value = 1.054462929432981
$ws.Range("E7").Value = 1.026870454065643
$ws.Range("F7").Value = 0.5940075024340167
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = 0.9660439600786086
$ws.Range("C8").Value = 1.214125379206839
$ws.Range("D8").Value = 2.475674515919399
$ws.Range("E8").Value = 1.573427632882872
$ws.Range("F8").Value = 1.360485317722336
$ws.Range("G8").Value = 6

$ws.Range("B9").Value = 1.443592723702684
$ws.Range("C9").Value = 1.443592723702684
$ws.Range("D9").Value = 2.873932502829098
$ws.Range("E9").Value = 1.695267678813319
$ws.Range("F9").Value = 1.088558141007014
$ws.Range("G9").Value = 3

# Add new row 10 (Q8) - copy formatting from A9 (style s="1") then set the new value
$ws.Range("A9").Copy($ws.Range("A10"))
$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = 0.4350014876132097
$ws.Range("C10").Value = 0.4350014876132097
$ws.Range("D10").Value = 0.1892262942257054
$ws.Range("E10").Value = 0.4350014876132097
$ws.Range("G10").Value = 1
